$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency market snapshot refresh (prices + 1h volume %, and two row swaps
# reflecting updated coinranking.com rank order). Each entry is (CellAddress, NewText).

# B/C (coin name / link) columns never look numeric, so a plain .Value assignment is safe.
$textCells = @(
    @{ Addr = 'B29'; Text = 'PancakeSwap' }
    @{ Addr = 'C29'; Text = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake' }
    @{ Addr = 'B30'; Text = 'Toncoin' }
    @{ Addr = 'C30'; Text = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' }
    @{ Addr = 'B42'; Text = 'Algorand' }
    @{ Addr = 'C42'; Text = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' }
    @{ Addr = 'B43'; Text = 'THORChain' }
    @{ Addr = 'C43'; Text = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune' }
    @{ Addr = 'B44'; Text = 'MultiversX' }
    @{ Addr = 'C44'; Text = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld' }
)

foreach ($item in $textCells) {
    $ws.Range($item.Addr).Value = $item.Text
}

# D/E (price / volume%) columns hold numeric-looking strings ("254.54", "  +3.61%  ").
# Assigning .Value directly lets Excel coerce them to real numbers, which would also
# change the cell from inline-string to numeric type. To keep them as literal text (as
# they were originally authored) without touching the cell style, write each as a
# quoted-literal formula, then Copy + PasteSpecial(xlPasteValues) to bake it down to a
# plain text value and drop the formula.
$xlPasteValues = -4163
$textLiteralCells = @(
    @{ Addr = 'D2'; Text = '43.036.35' }
    @{ Addr = 'E2'; Text = '  +0.24%  ' }
    @{ Addr = 'D3'; Text = '2.220.34' }
    @{ Addr = 'E3'; Text = '  -1.19%  ' }
    @{ Addr = 'E4'; Text = '  -0.24%  ' }
    @{ Addr = 'D5'; Text = '254.54' }
    @{ Addr = 'E5'; Text = '  +3.61%  ' }
    @{ Addr = 'D6'; Text = '0.615' }
    @{ Addr = 'E6'; Text = '  -0.29%  ' }
    @{ Addr = 'D7'; Text = '76.24' }
    @{ Addr = 'E7'; Text = '  +0.20%  ' }
    @{ Addr = 'E8'; Text = '  +0.01%  ' }
    @{ Addr = 'D9'; Text = '0.593' }
    @{ Addr = 'E9'; Text = '  -3.97%  ' }
    @{ Addr = 'D10'; Text = '41.42' }
    @{ Addr = 'E10'; Text = '  +0.73%  ' }
    @{ Addr = 'D11'; Text = '0.0916' }
    @{ Addr = 'E11'; Text = '  -2.14%  ' }
    @{ Addr = 'D12'; Text = '6.91' }
    @{ Addr = 'E12'; Text = '  -1.08%  ' }
    @{ Addr = 'E13'; Text = '  +0.30%  ' }
    @{ Addr = 'D14'; Text = '2.552.64' }
    @{ Addr = 'E14'; Text = '  -0.34%  ' }
    @{ Addr = 'D15'; Text = '14.36' }
    @{ Addr = 'E15'; Text = '  -1.98%  ' }
    @{ Addr = 'D16'; Text = '2.221.73' }
    @{ Addr = 'D17'; Text = '0.782' }
    @{ Addr = 'E17'; Text = '  -3.61%  ' }
    @{ Addr = 'D18'; Text = '42.928.05' }
    @{ Addr = 'E18'; Text = '  +0.18%  ' }
    @{ Addr = 'D19'; Text = '0.0000103' }
    @{ Addr = 'E19'; Text = '  -2.25%  ' }
    @{ Addr = 'D20'; Text = '71.18' }
    @{ Addr = 'E20'; Text = '  -0.23%  ' }
    @{ Addr = 'D21'; Text = '5.95' }
    @{ Addr = 'E21'; Text = '  -0.96%  ' }
    @{ Addr = 'D22'; Text = '229.64' }
    @{ Addr = 'E22'; Text = '  -0.67%  ' }
    @{ Addr = 'D23'; Text = '2.19' }
    @{ Addr = 'E23'; Text = '  -1.56%  ' }
    @{ Addr = 'D24'; Text = '9.37' }
    @{ Addr = 'E24'; Text = '  -7.81%  ' }
    @{ Addr = 'E25'; Text = '  +0.06%  ' }
    @{ Addr = 'D26'; Text = '10.61' }
    @{ Addr = 'E26'; Text = '  -2.80%  ' }
    @{ Addr = 'D27'; Text = '3.37' }
    @{ Addr = 'E27'; Text = '  -1.05%  ' }
    @{ Addr = 'D28'; Text = '38.68' }
    @{ Addr = 'E28'; Text = '  +1.40%  ' }
    @{ Addr = 'D29'; Text = '2.19' }
    @{ Addr = 'E29'; Text = '  -2.79%  ' }
    @{ Addr = 'D30'; Text = '2.21' }
    @{ Addr = 'E30'; Text = '  +2.84%  ' }
    @{ Addr = 'D31'; Text = '173.25' }
    @{ Addr = 'E31'; Text = '  +0.06%  ' }
    @{ Addr = 'D32'; Text = '20.26' }
    @{ Addr = 'E32'; Text = '  -0.31%  ' }
    @{ Addr = 'D33'; Text = '0.0846' }
    @{ Addr = 'E33'; Text = '  +6.53%  ' }
    @{ Addr = 'D34'; Text = '5.23' }
    @{ Addr = 'E34'; Text = '  -2.58%  ' }
    @{ Addr = 'E35'; Text = '  -1.35%  ' }
    @{ Addr = 'E36'; Text = '  -0.87%  ' }
    @{ Addr = 'D37'; Text = '0.0351' }
    @{ Addr = 'E37'; Text = '  +7.22%  ' }
    @{ Addr = 'D38'; Text = '4.31' }
    @{ Addr = 'E38'; Text = '  -0.96%  ' }
    @{ Addr = 'D39'; Text = '12.52' }
    @{ Addr = 'E39'; Text = '  -2.51%  ' }
    @{ Addr = 'D40'; Text = '2.11' }
    @{ Addr = 'E40'; Text = '  -1.78%  ' }
    @{ Addr = 'E41'; Text = '  +17.38%  ' }
    @{ Addr = 'D42'; Text = '0.198' }
    @{ Addr = 'E42'; Text = '  -3.13%  ' }
    @{ Addr = 'D43'; Text = '5.27' }
    @{ Addr = 'E43'; Text = '  -5.46%  ' }
    @{ Addr = 'D44'; Text = '60.10' }
    @{ Addr = 'E44'; Text = '  +0.36%  ' }
    @{ Addr = 'D45'; Text = '102.03' }
    @{ Addr = 'E45'; Text = '  -4.45%  ' }
    @{ Addr = 'D46'; Text = '8.34' }
    @{ Addr = 'E46'; Text = '  -3.90%  ' }
    @{ Addr = 'D47'; Text = '0.0978' }
    @{ Addr = 'E47'; Text = '  -1.27%  ' }
    @{ Addr = 'D48'; Text = '0.460' }
    @{ Addr = 'E48'; Text = '  +2.43%  ' }
    @{ Addr = 'D49'; Text = '1.11' }
    @{ Addr = 'E49'; Text = '  +0.12%  ' }
    @{ Addr = 'D50'; Text = '1.14' }
    @{ Addr = 'E50'; Text = '  -1.04%  ' }
    @{ Addr = 'D51'; Text = '2.446.72' }
    @{ Addr = 'E51'; Text = '  -0.01%  ' }
)

foreach ($item in $textLiteralCells) {
    $cell = $ws.Range($item.Addr)
    $cell.Formula = '="' + $item.Text + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial($xlPasteValues) | Out-Null
}

$excel.CutCopyMode = 0
